$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.267.11"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "3.491.19"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.38"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.91"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").Value = "3.484.71"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.93"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "4.091.73"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.69"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "67.184.92"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000174"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "3.490.14"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.03"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.96"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.85"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.90"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.545"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.81"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.92"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000122"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.54"
$ws.Range("E31").Value = "  +5.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.94"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.23"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  +14.76%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.897"
$ws.Range("E40").Value = "  +2.84%  "
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.66"
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.755.61"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.46"
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0702"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.74"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.85"
$ws.Range("E48").Value = "  -4.00%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "316.84"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("E51").Value = "  -1.09%  "
